$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "auto"
$ws.Range("B1").Value = "l2"
$ws.Range("C1").Value = "Lâmpada"
$ws.Range("D1").Value = 100
$ws.Range("E1").Value = $false
